$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 6).Value = 56
$ws.Cells.Item(6, 7).Value = 8
$ws.Cells.Item(15, 6).Value = 219
$ws.Cells.Item(23, 6).Value = 535
$ws.Cells.Item(23, 7).Value = 4
$ws.Cells.Item(35, 6).Value = 2
$ws.Cells.Item(236, 6).Value = 2176
$ws.Cells.Item(236, 7).Value = 231
$ws.Cells.Item(240, 6).Value = 41436
$ws.Cells.Item(240, 7).Value = 494
$ws.Cells.Item(241, 6).Value = 89780
$ws.Cells.Item(244, 6).Value = 5283
$ws.Cells.Item(244, 7).Value = 98
$ws.Cells.Item(245, 6).Value = 3606
$ws.Cells.Item(245, 7).Value = 82
$ws.Cells.Item(246, 6).Value = 2115
$ws.Cells.Item(246, 7).Value = 140
$ws.Cells.Item(250, 6).Value = 17658
$ws.Cells.Item(250, 7).Value = 1154
$ws.Cells.Item(257, 6).Value = 5642
$ws.Cells.Item(257, 7).Value = 286
$ws.Cells.Item(258, 6).Value = 3919
$ws.Cells.Item(258, 7).Value = 224
$ws.Cells.Item(259, 6).Value = 6546
$ws.Cells.Item(259, 7).Value = 442
$ws.Cells.Item(260, 6).Value = 12193
$ws.Cells.Item(260, 7).Value = 770
$ws.Cells.Item(261, 6).Value = 18110
$ws.Cells.Item(261, 7).Value = 611
$ws.Cells.Item(262, 6).Value = 9285
$ws.Cells.Item(262, 7).Value = 322
$ws.Cells.Item(264, 6).Value = 43207
$ws.Cells.Item(264, 7).Value = 885
$ws.Cells.Item(265, 6).Value = 18446
$ws.Cells.Item(265, 7).Value = 894
$ws.Cells.Item(266, 6).Value = 14726
$ws.Cells.Item(266, 7).Value = 733
$ws.Cells.Item(267, 6).Value = 15513
$ws.Cells.Item(267, 7).Value = 870
$ws.Cells.Item(268, 6).Value = 17539
$ws.Cells.Item(268, 7).Value = 787
$ws.Cells.Item(269, 6).Value = 9705
$ws.Cells.Item(269, 7).Value = 433
$ws.Cells.Item(270, 6).Value = 2780
$ws.Cells.Item(270, 7).Value = 180
$ws.Cells.Item(272, 6).Value = 30785
$ws.Cells.Item(272, 7).Value = 1662
$ws.Cells.Item(273, 6).Value = 31709
$ws.Cells.Item(273, 7).Value = 1662
$ws.Cells.Item(274, 6).Value = 28090
$ws.Cells.Item(274, 7).Value = 1277
$ws.Cells.Item(275, 6).Value = 30330
$ws.Cells.Item(277, 6).Value = 3393
$ws.Cells.Item(278, 6).Value = 30542
$ws.Cells.Item(279, 6).Value = 42700
$ws.Cells.Item(279, 7).Value = 3036
$ws.Cells.Item(280, 6).Value = 34725
$ws.Cells.Item(280, 7).Value = 2310
$ws.Cells.Item(281, 6).Value = 46078
$ws.Cells.Item(286, 6).Value = 55193
$ws.Cells.Item(287, 6).Value = 58840
$ws.Cells.Item(287, 7).Value = 3717
$ws.Cells.Item(288, 6).Value = 59278
$ws.Cells.Item(289, 6).Value = 62984
$ws.Cells.Item(290, 6).Value = 17572
$ws.Cells.Item(294, 6).Value = 93811
$ws.Cells.Item(294, 7).Value = 4942
$ws.Cells.Item(300, 6).Value = 72571
$ws.Cells.Item(300, 7).Value = 6981
$ws.Cells.Item(301, 6).Value = 72136
$ws.Cells.Item(301, 7).Value = 5676
$ws.Cells.Item(302, 6).Value = 78622
$ws.Cells.Item(302, 7).Value = 5656
$ws.Cells.Item(307, 6).Value = 75826
$ws.Cells.Item(307, 7).Value = 6395
$ws.Cells.Item(308, 6).Value = 15395
$ws.Cells.Item(309, 6).Value = 77913
$ws.Cells.Item(314, 6).Value = 64302
$ws.Cells.Item(314, 7).Value = 3148
$ws.Cells.Item(315, 6).Value = 56287
$ws.Cells.Item(315, 7).Value = 2627
$ws.Cells.Item(316, 6).Value = 50749
$ws.Cells.Item(317, 6).Value = 63737
$ws.Cells.Item(321, 6).Value = 89345
$ws.Cells.Item(321, 7).Value = 2654
$ws.Cells.Item(322, 6).Value = 109547
$ws.Cells.Item(323, 6).Value = 216834
$ws.Cells.Item(324, 6).Value = 241036
$ws.Cells.Item(325, 6).Value = 766038
$ws.Cells.Item(325, 7).Value = 6455
$ws.Cells.Item(326, 6).Value = 419543
$ws.Cells.Item(327, 6).Value = 225126
$ws.Cells.Item(327, 7).Value = 2721
$ws.Cells.Item(328, 6).Value = 180664
$ws.Cells.Item(328, 7).Value = 2665
$ws.Cells.Item(329, 6).Value = 82996
$ws.Cells.Item(329, 7).Value = 1757
$ws.Cells.Item(330, 6).Value = 72545
$ws.Cells.Item(330, 7).Value = 2083
$ws.Cells.Item(331, 6).Value = 154995
$ws.Cells.Item(332, 6).Value = 457203
$ws.Cells.Item(333, 6).Value = 271706
$ws.Cells.Item(333, 7).Value = 2942
$ws.Cells.Item(334, 6).Value = 196773
$ws.Cells.Item(335, 6).Value = 130851
$ws.Cells.Item(335, 7).Value = 2998
$ws.Cells.Item(336, 6).Value = 102552
$ws.Cells.Item(336, 7).Value = 3339
$ws.Cells.Item(337, 6).Value = 103557
$ws.Cells.Item(337, 7).Value = 2890
$ws.Cells.Item(338, 6).Value = 228385
$ws.Cells.Item(338, 7).Value = 3202
$ws.Cells.Item(339, 6).Value = 661287
$ws.Cells.Item(341, 6).Value = 291759
$ws.Cells.Item(342, 6).Value = 178654
$ws.Cells.Item(342, 7).Value = 3038
$ws.Cells.Item(344, 6).Value = 135478
$ws.Cells.Item(346, 6).Value = 676087
$ws.Cells.Item(346, 7).Value = 4837
$ws.Cells.Item(347, 6).Value = 343700
$ws.Cells.Item(348, 6).Value = 232622
$ws.Cells.Item(348, 7).Value = 3240
$ws.Cells.Item(349, 6).Value = 159176
$ws.Cells.Item(349, 7).Value = 2754
$ws.Cells.Item(350, 6).Value = 127248
$ws.Cells.Item(350, 7).Value = 2790
$ws.Cells.Item(351, 6).Value = 150929
$ws.Cells.Item(351, 7).Value = 2832
$ws.Cells.Item(352, 6).Value = 307405
$ws.Cells.Item(352, 7).Value = 3542
$ws.Cells.Item(353, 6).Value = 725584
$ws.Cells.Item(353, 7).Value = 5304
$ws.Cells.Item(355, 6).Value = 222058
$ws.Cells.Item(355, 7).Value = 3448
$ws.Cells.Item(356, 6).Value = 159963
$ws.Cells.Item(357, 6).Value = 138591
$ws.Cells.Item(357, 7).Value = 3029
$ws.Cells.Item(358, 6).Value = 158739
$ws.Cells.Item(358, 7).Value = 2613
$ws.Cells.Item(359, 6).Value = 321217
$ws.Cells.Item(359, 7).Value = 3342
$ws.Cells.Item(360, 6).Value = 751668
$ws.Cells.Item(360, 7).Value = 5141
$ws.Cells.Item(361, 6).Value = 332572
$ws.Cells.Item(361, 7).Value = 2622
$ws.Cells.Item(362, 6).Value = 229127
$ws.Cells.Item(362, 7).Value = 3182
$ws.Cells.Item(363, 6).Value = 188458
$ws.Cells.Item(363, 7).Value = 2756
$ws.Cells.Item(364, 6).Value = 168445
$ws.Cells.Item(364, 7).Value = 2480
$ws.Cells.Item(365, 6).Value = 184554
$ws.Cells.Item(365, 7).Value = 2398
$ws.Cells.Item(366, 6).Value = 340159
$ws.Cells.Item(366, 7).Value = 2850
$ws.Cells.Item(367, 6).Value = 766203
$ws.Cells.Item(368, 6).Value = 346351
$ws.Cells.Item(368, 7).Value = 2300
$ws.Cells.Item(369, 6).Value = 234542
$ws.Cells.Item(369, 7).Value = 2600
$ws.Cells.Item(370, 6).Value = 180653
$ws.Cells.Item(370, 7).Value = 2046
$ws.Cells.Item(371, 6).Value = 160080
$ws.Cells.Item(372, 6).Value = 179774
$ws.Cells.Item(372, 7).Value = 1869
$ws.Cells.Item(373, 6).Value = 349537
$ws.Cells.Item(373, 7).Value = 2370
$ws.Cells.Item(374, 6).Value = 772911
$ws.Cells.Item(374, 7).Value = 3423
$ws.Cells.Item(375, 6).Value = 351597
$ws.Cells.Item(375, 7).Value = 1857
$ws.Cells.Item(376, 6).Value = 220919
$ws.Cells.Item(377, 6).Value = 176820
$ws.Cells.Item(377, 7).Value = 1813
$ws.Cells.Item(378, 6).Value = 157605
$ws.Cells.Item(378, 7).Value = 1563
$ws.Cells.Item(379, 6).Value = 179746
$ws.Cells.Item(379, 7).Value = 1611
$ws.Cells.Item(380, 6).Value = 344765
$ws.Cells.Item(380, 7).Value = 2022
$ws.Cells.Item(381, 6).Value = 744261
$ws.Cells.Item(381, 7).Value = 2684
$ws.Cells.Item(382, 6).Value = 357646
$ws.Cells.Item(383, 6).Value = 220770
$ws.Cells.Item(384, 6).Value = 171986
$ws.Cells.Item(384, 7).Value = 1513
$ws.Cells.Item(385, 6).Value = 150852
$ws.Cells.Item(385, 7).Value = 1409
$ws.Cells.Item(386, 6).Value = 182480
$ws.Cells.Item(386, 7).Value = 1360
$ws.Cells.Item(388, 6).Value = 728552
$ws.Cells.Item(389, 6).Value = 353116
$ws.Cells.Item(390, 6).Value = 220123
$ws.Cells.Item(390, 7).Value = 1475
$ws.Cells.Item(391, 6).Value = 176838
$ws.Cells.Item(391, 7).Value = 1208
$ws.Cells.Item(392, 6).Value = 220410
$ws.Cells.Item(392, 7).Value = 1217
$ws.Cells.Item(393, 6).Value = 299587
$ws.Cells.Item(393, 7).Value = 1199
$ws.Cells.Item(394, 6).Value = 162415
$ws.Cells.Item(395, 6).Value = 738885
$ws.Cells.Item(396, 6).Value = 164183
$ws.Cells.Item(396, 7).Value = 548
$ws.Cells.Item(397, 6).Value = 106344
$ws.Cells.Item(397, 7).Value = 629
$ws.Cells.Item(398, 6).Value = 291580
$ws.Cells.Item(398, 7).Value = 1447
$ws.Cells.Item(399, 6).Value = 195274
$ws.Cells.Item(400, 6).Value = 145532
$ws.Cells.Item(400, 7).Value = 736
$ws.Cells.Item(401, 6).Value = 264469
$ws.Cells.Item(401, 7).Value = 921
$ws.Cells.Item(402, 6).Value = 695722
$ws.Cells.Item(402, 7).Value = 1340
$ws.Cells.Item(403, 6).Value = 344587
$ws.Cells.Item(403, 7).Value = 722
$ws.Cells.Item(404, 6).Value = 217928
$ws.Cells.Item(404, 7).Value = 885
$ws.Cells.Item(405, 6).Value = 165004
$ws.Cells.Item(405, 7).Value = 664

Write-Output "Applied 221 cell updates"
